$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data table changed (corrected/re-simulated values) and one new
# row (row 53) was appended, so the whole A2:E53 block is rewritten in place.

# Row 53 is brand new: first copy the formatting (style) of the last existing
# data row (52) down into it, so A53 keeps the same date style as the rest of column A.
$ws.Range("A52").Copy($ws.Range("A53"))

$data = @(
    @(39400, 2007, -0.7015558851707349, 2008, 1.245556224328537),
    @(39583, 2008, 1.247512323997491, 2009, 2.82953744009995),
    @(39765, 2008, 0.3590181115727287, 2009, 2.957479223435766),
    @(39948, 2009, 0.5917823527752164, 2010, -1.985049937499994),
    @(40130, 2009, -0.01587181126743165, 2010, 3.835144370775678),
    @(40310, 2010, 1.985742476411234, 2011, 7.819356632099961),
    @(40494, 2010, -0.02256889165885845, 2011, -0.6757980944263386),
    @(40676, 2011, -0.7518797681958955, 2012, -1.590425574400001),
    @(40862, 2011, 0.09611428386597787, 2012, -1.14533080410083),
    @(41044, 2012, -0.5765930039052902, 2013, -1.590425574400001),
    @(41228, 2012, -0.1827723404408288, 2013, -0.9712385417624092),
    @(41409, 2013, -0.07642926654481963, 2014, -0.3994003998999851),
    @(41592, 2013, -0.001350220946483294, 2014, 0.6008487920565297),
    @(41774, 2014, 0.9274109147535681, 2015, 2.829537440099972),
    @(41957, 2014, -0.075754880139145, 2015, -1.589980933184099),
    @(42137, 2015, -0.5259734324162268, 2016, -1.194610791900008),
    @(42321, 2015, -0.5761528471665445, 2016, 0.4501721032283079),
    @(42503, 2016, 0.07456754038981384, 2017, 0.4006004000999708),
    @(42689, 2016, -0.2011999787958185, 2017, -0.1501751374934579),
    @(42867, 2017, -0.07666472728168339, 2018, -1.590425574400012),
    @(43053, 2017, 0.1213692818849532, 2018, 1.204188598110267),
    @(43145, 2018, 0.1986438914956645, 2019, 0),
    @(43235, 2018, -0.2004754673795017, 2019, -0.3994003998999962),
    @(43326, 2018, 0.1743342151774741, 2019, 1.483403742553846),
    @(43418, 2018, 0.1493219406571766, 2019, 2.372078088364726),
    @(43510, 2019, -0.7530239469328848, 2020, -2.378486270400004),
    @(43600, 2019, -0.7283174404322912, 2020, -2.378486270400004),
    @(43691, 2019, -0.2044553505917923, 2020, 0.8757395670492052),
    @(43783, 2019, -0.4278219446121612, 2020, -2.378564786744741),
    @(43875, 2020, 0.6687116611293176, 2021, 2.42168652960002),
    @(43966, 2020, -0.03096525636255842, 2021, 1.205410808099971),
    @(44068, 2020, -0.4412356890029168, 2021, -1.287161643752965),
    @(44159, 2020, -1.026566979837418, 2021, -3.329887238705409),
    @(44251, 2021, -0.9187270897784594, 2022, -2.100316115087164),
    @(44341, 2021, -0.04074803603358879, 2022, -0.2354831990173722),
    @(44432, 2021, -0.04976849661378902, 2022, -6.367504067266704),
    @(44525, 2021, 0.3179894933462268, 2022, 0.09136665742732752),
    @(44617, 2022, 0.7185676380534911, 2023, 1.072043735370976),
    @(44706, 2022, 1.296559640836992, 2023, 0.439921111559638),
    @(44798, 2022, 0.4784173072842179, 2023, -0.07349332085010429),
    @(44890, 2022, 0.463604920919658, 2023, -1.086119702215593),
    @(44981, 2023, 0.3516957773510798, 2024, -0.7234699379142895),
    @(45071, 2023, 0.8143067496459322, 2024, -0.05551881767973388),
    @(45163, 2023, 0.4279153732809959, 2024, 0.09327584436471525),
    @(45254, 2023, 0.621639092134818, 2024, -0.185520708596465),
    @(45345, 2024, -0.4267958602204081, 2025, 0.00005213214835375624),
    @(45436, 2024, -1.077932696718564, 2025, -0.02296885644381685),
    @(45534, 2024, -0.9378224616154895, 2025, 0.3335355268917262),
    @(45618, 2024, -0.6768900623516982, 2025, 0.9196327481985289),
    @(45713, 2025, 2.173733407444822, 2026, 1.89918335857393),
    @(45800, 2025, 0.4743170952486997, 2026, -0.324019587165425),
    @(45891, 2025, 1.136769786738334, 2026, 0.2972123973886909)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Host "done"
